# Add a "2022-Q3" sheet of fund-holding data ahead of the existing
# "2021-Q3" sheet, and record the new quarter's totals as the newest
# (top) row on the "总计" (totals) sheet.
#
# Starting layout:  总计 (sheetId 1), 2021-Q3 (sheetId 2)
# Target layout:    总计 (sheetId 1), 2022-Q3 (sheetId 2), 2021-Q3 (sheetId 3)
#
# To get the "2021-Q3" data onto a *new* sheetId 3 (rather than leaving it on
# sheetId 2 and giving the new data sheetId 3), duplicate the existing
# "2021-Q3" sheet first, rename the original to "2022-Q3" and overwrite its
# contents with the new data, then rename the untouched duplicate back to
# "2021-Q3". That reproduces both the id/name remap *and* keeps the old
# sheet's data byte-for-byte.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2021-Q3")

# Duplicate "2021-Q3" immediately after itself, then shuffle names so the
# ORIGINAL sheet object (still sheetId 2) becomes "2022-Q3", and the COPY
# (new sheetId 3) keeps the "2021-Q3" name/data.
$q3.Copy($null, $q3)
$q3.Name = "2022-Q3"
$q3copy = $wb.Worksheets.Item(3)
$q3copy.Name = "2021-Q3"

# $q3 is now the "2022-Q3" sheet; wipe its old 2021-Q3 fund rows.
$q3.UsedRange.Clear()

# ---- header row (B1:H1), styled like the bold/centered/bordered header
# cells already used elsewhere in the workbook (copy format only) ----
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# ---- data rows ----
# Columns D,E,F,G (and B) hold text that looks numeric (fund codes/percents),
# so force text entry via NumberFormat "@" first to avoid Excel silently
# re-typing them as numbers (and dropping leading zeros / changing 0.78 ->
# 0.78000000000000003). G6 is a genuine number (0), and column H is numeric.
$data = @(
    @(0, "161224", "国投瑞银新丝路灵活配置混合（LOF）", "0.78", "93.98", "4.41", "0.0344", 8),
    @(1, "012432", "国投瑞银安泰混合C", "1.00", "32.06", "2.77", "0.0277", 5),
    @(2, "012019", "国投瑞银安泽混合A", "0.49", "32.69", "2.43", "0.0119", 5),
    @(3, "012020", "国投瑞银安泽混合C", "0.10", "32.69", "2.43", "0.0024", 5),
    @(4, "012431", "国投瑞银安泰混合A", "0.00", "32.06", "2.77", $null, 5)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = 2 + $r
    $rowvals = $data[$r]

    $q3.Cells.Item($row, 1).Value = $rowvals[0]

    $cB = $q3.Cells.Item($row, 2)
    $cB.NumberFormat = "@"
    $cB.Value = $rowvals[1]

    $q3.Cells.Item($row, 3).Value = $rowvals[2]

    $cD = $q3.Cells.Item($row, 4)
    $cD.NumberFormat = "@"
    $cD.Value = $rowvals[3]

    $cE = $q3.Cells.Item($row, 5)
    $cE.NumberFormat = "@"
    $cE.Value = $rowvals[4]

    $cF = $q3.Cells.Item($row, 6)
    $cF.NumberFormat = "@"
    $cF.Value = $rowvals[5]

    $cG = $q3.Cells.Item($row, 7)
    if ($rowvals[6] -eq $null) {
        $cG.Value = 0
    } else {
        $cG.NumberFormat = "@"
        $cG.Value = $rowvals[6]
    }

    $q3.Cells.Item($row, 8).Value = $rowvals[7]
}

# Clear the temporary "@" number-format back to General (the stored values
# stay text even after the display format changes), matching the unstyled
# data cells elsewhere in the workbook.
$total.Range("C2").Copy()
$q3.Range("B2:G6").PasteSpecial(-4122)

# First column (index 0..4) uses the same bold/centered/bordered style as
# the "总计" sheet's index column.
$total.Range("A2").Copy()
$q3.Range("A2:A6").PasteSpecial(-4122)

# ---- update "总计": new quarter's row becomes row 2, old 2021-Q3 row
# shifts down to row 3 (values/style preserved) ----
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2021-Q3"
$total.Cells.Item(3, 3).Value = 1
$total.Cells.Item(3, 4).Value = 0.08
$total.Range("A2").Copy()
$total.Cells.Item(3, 1).PasteSpecial(-4122)

$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.08
